# Update cryptocurrency price/volume table to latest scrape.
# Row 12/13 (Solana/BinanceUSD) and row 27/28 (LidoDAOToken/EthereumClassic)
# swapped ranking order, so B/C/D/E are all rewritten for those rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''27.965.91'
$ws.Range("E2").Value = '  -0.83%  '

# Row 3
$ws.Range("D3").Value = '''1.763.76'
$ws.Range("E3").Value = '  -3.08%  '

# Row 4
$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  +0.77%  '

# Row 5
$ws.Range("D5").Value = '''338.92'
$ws.Range("E5").Value = '  -0.15%  '

# Row 6
$ws.Range("D6").Value = '''0.9977'
$ws.Range("E6").Value = '  +0.37%  '

# Row 7
$ws.Range("D7").Value = '''0.3768'
$ws.Range("E7").Value = '  -4.22%  '

# Row 8
$ws.Range("D8").Value = '''0.3369'
$ws.Range("E8").Value = '  -3.68%  '

# Row 9
$ws.Range("D9").Value = '''46.16'
$ws.Range("E9").Value = '  -4.68%  '

# Row 10
$ws.Range("D10").Value = '''1.131'
$ws.Range("E10").Value = '  -6.04%  '

# Row 11
$ws.Range("D11").Value = '''0.07212'
$ws.Range("E11").Value = '  -5.23%  '

# Row 12
$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D12").Value = '''1.002'
$ws.Range("E12").Value = '  +0.88%  '

# Row 13
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").Value = '''22.62'
$ws.Range("E13").Value = '  +1.49%  '

# Row 14
$ws.Range("D14").Value = '''6.202'
$ws.Range("E14").Value = '  -5.46%  '

# Row 15
$ws.Range("D15").Value = '''7.175'
$ws.Range("E15").Value = '  -0.56%  '

# Row 16
$ws.Range("D16").Value = '''1.763.39'
$ws.Range("E16").Value = '  -2.80%  '

# Row 17
$ws.Range("D17").Value = '''0.00001059'
$ws.Range("E17").Value = '  -4.63%  '

# Row 18
$ws.Range("D18").Value = '''0.06580'
$ws.Range("E18").Value = '  -1.83%  '

# Row 19
$ws.Range("D19").Value = '''80.74'

# Row 20
$ws.Range("D20").Value = '''0.9986'
$ws.Range("E20").Value = '  +0.31%  '

# Row 21
$ws.Range("D21").Value = '''16.97'
$ws.Range("E21").Value = '  -5.29%  '

# Row 22
$ws.Range("D22").Value = '''6.272'
$ws.Range("E22").Value = '  -4.77%  '

# Row 23
$ws.Range("D23").Value = '''27.994.59'
$ws.Range("E23").Value = '  -0.66%  '

# Row 24
$ws.Range("D24").Value = '''11.71'
$ws.Range("E24").Value = '  -8.82%  '

# Row 25
$ws.Range("D25").Value = '''2.398'
$ws.Range("E25").Value = '  -0.47%  '

# Row 26
$ws.Range("D26").Value = '''153.43'
$ws.Range("E26").Value = '  -0.92%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '''19.83'
$ws.Range("E27").Value = '  -7.44%  '

# Row 28
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '''2.340'
$ws.Range("E28").Value = '  -9.55%  '

# Row 29
$ws.Range("D29").Value = '''1.286'
$ws.Range("E29").Value = '  -16.59%  '

# Row 30
$ws.Range("D30").Value = '''1.964.90'
$ws.Range("E30").Value = '  -2.75%  '

# Row 31
$ws.Range("D31").Value = '''131.30'
$ws.Range("E31").Value = '  -3.30%  '

# Row 32
$ws.Range("D32").Value = '''4.017'
$ws.Range("E32").Value = '  -0.45%  '

# Row 33
$ws.Range("D33").Value = '''5.848'
$ws.Range("E33").Value = '  -5.71%  '

# Row 34
$ws.Range("D34").Value = '''0.08786'
$ws.Range("E34").Value = '  -0.76%  '

# Row 35
$ws.Range("D35").Value = '''12.27'
$ws.Range("E35").Value = '  -8.30%  '

# Row 36
$ws.Range("D36").Value = '''0.02343'
$ws.Range("E36").Value = '  -3.97%  '

# Row 37
$ws.Range("D37").Value = '''0.6603'
$ws.Range("E37").Value = '  -5.44%  '

# Row 38
$ws.Range("D38").Value = '''0.06222'
$ws.Range("E38").Value = '  -5.26%  '

# Row 39
$ws.Range("D39").Value = '''5.157'
$ws.Range("E39").Value = '  -7.12%  '

# Row 40
$ws.Range("D40").Value = '''0.2114'
$ws.Range("E40").Value = '  -5.49%  '

# Row 41
$ws.Range("D41").Value = '''1.208'
$ws.Range("E41").Value = '  -4.85%  '

# Row 42
$ws.Range("D42").Value = '''1.452'
$ws.Range("E42").Value = '  -10.02%  '

# Row 43
$ws.Range("D43").Value = '''8.046'
$ws.Range("E43").Value = '  -6.29%  '

# Row 44
$ws.Range("D44").Value = '''0.9982'
$ws.Range("E44").Value = '  +0.38%  '

# Row 45
$ws.Range("D45").Value = '''13.77'
$ws.Range("E45").Value = '  -5.99%  '

# Row 46
$ws.Range("D46").Value = '''3.835'
$ws.Range("E46").Value = '  -1.06%  '

# Row 47
$ws.Range("D47").Value = '''0.6059'
$ws.Range("E47").Value = '  -7.30%  '

# Row 48
$ws.Range("D48").Value = '''130.29'
$ws.Range("E48").Value = '  -1.79%  '

# Row 49
$ws.Range("D49").Value = '''2.015'
$ws.Range("E49").Value = '  -7.33%  '

# Row 50
$ws.Range("D50").Value = '''0.07249'
$ws.Range("E50").Value = '  +0.34%  '

# Row 51
$ws.Range("D51").Value = '''1.182'
$ws.Range("E51").Value = '  +1.67%  '
